# "2017 ppt fixes" - update the title-slide subtitle from
# "FME 2016 Training" to "FME " + "2017" (typed as two runs, with the
# trailing " Training" removed), matching the authored OOXML diff.
#
# Note: the diff also swaps which embedded-font entry
# (Proxima Nova / Open Sans) in ppt/presentation.xml's
# <p:embeddedFontLst> claims each already-embedded font binary
# (rId23-26 vs rId27-30 are untouched). That is pure embedded-font
# metadata that PowerPoint manages internally when it embeds fonts on
# save; it is not reachable through the Presentation/Shape/TextRange
# object model (no Fonts.Item(...).Name-style setter actually mutates
# it here), so it is intentionally left alone rather than risk
# corrupting the embedded font relationships.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape 2 on slide 1 is the "subTitle" placeholder shape containing the
# "FME 2016 Training" text.
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$prefix = "FME "          # kept as-is -> becomes its own run
$oldYear = "2016"          # retyped -> "2017", splitting off a new run
$suffix = " Training"     # deleted entirely

$yearStart = $prefix.Length + 1
$tailStart = $yearStart + $oldYear.Length

# 1) Delete the trailing " Training" text.
if ($tr.Length -ge $tailStart) {
    $tail = $tr.Characters($tailStart, $tr.Length - ($tailStart - 1))
    $tail.Text = ""
}

# 2) Retype "2016" as "2017". Editing just this sub-range leaves "FME "
#    and "2017" as two separate runs, matching the target markup.
$year = $tr.Characters($yearStart, $oldYear.Length)
$year.Text = "2017"
